# Update the "time_taken" column (F) on the "data" sheet with refreshed timestamps
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:21:52.638019",
    "2021-10-05 14:21:52.638027",
    "2021-10-05 14:21:52.638031",
    "2021-10-05 14:21:52.638033",
    "2021-10-05 14:21:52.638036",
    "2021-10-05 14:21:52.638039",
    "2021-10-05 14:21:52.638042",
    "2021-10-05 14:21:52.638044",
    "2021-10-05 14:21:52.638047",
    "2021-10-05 14:21:52.638049",
    "2021-10-05 14:21:52.638052",
    "2021-10-05 14:21:52.638054",
    "2021-10-05 14:21:52.638057",
    "2021-10-05 14:21:52.638059",
    "2021-10-05 14:21:52.638062"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Add a new "metadata" worksheet positioned right after the "data" sheet
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row: copy the bold/centered header formatting from the "data" sheet, then set labels
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# First (and only) data row
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$metaSheet.Range("A2").Value = 0

$metaSheet.Range("B2").Value = "Ocular and oculo-cutaneous albinism"
$metaSheet.Range("C2").Value = 128
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.21"
$metaSheet.Range("D2").Style = $dataSheet.Range("D2").Style
$metaSheet.Range("E2").Value = "2019-06-20T15:13:41.618714Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:21:52.634285"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/128/?format=json"
